$d = $word.ActiveDocument

# --- Title ---
$d.Content.Find.Execute("Artificial Intelligence: Navigating the Ethical Conundrum", $true, $false, $false, $false, $false, $true, 1, $false, "The Fascinating Realm of Cells: Microscopic Worlds", 2)

# --- Author name ---
$d.Content.Find.Execute("Kevin Martin", $true, $false, $false, $false, $false, $true, 1, $false, "Dr. Alecia Marshall", 2)

# --- Email address parts ---
$d.Content.Find.Execute("kevin", $true, $true, $false, $false, $false, $true, 1, $false, "alecia", 2)
$d.Content.Find.Execute("martin@abcxyz", $true, $false, $false, $false, $false, $true, 1, $false, "marshall@educator", 2)
$d.Content.Find.Execute("com", $true, $true, $false, $false, $false, $true, 1, $false, "org", 2)

# --- Body paragraph 1 ---
$d.Content.Find.Execute("As artificial intelligence (AI) strides forward with remarkable progress, its pervasive integration into our lives brings ethical questions to the forefront", $true, $false, $false, $false, $false, $true, 1, $false, "The realm of cells is vast, teeming with countless microscopic worlds that hold the secrets of life", 2)

$d.Content.Find.Execute(" The intricate dance between human autonomy and AI-driven decision-making demands careful consideration", $true, $false, $false, $false, $false, $true, 1, $false, " These tiny, intricate structures are the foundation of all living things, from the towering redwood to the minuscule bacterium", 2)

$d.Content.Find.Execute(" We must unravel the potential apprehensions and dilemmas while exploring the positive applications of AI that enhance human lives", $true, $false, $false, $false, $false, $true, 1, $false, " Within these minute boundaries, cells perform complex symphonies of biochemical reactions, carrying out functions that sustain life and support growth", 2)

$d.Content.Find.Execute(" Striking a balance between progress and responsibility becomes crucial, ensuring AI's impact aligns with ethical standards and societal values", $true, $false, $false, $false, $false, $true, 1, $false, " Exploring the world of cells is a captivating journey into the fundamental mechanisms of biology, revealing the building blocks of life and shedding light on the mysteries of our own existence", 2)

$d.Content.Find.Execute("The transformative power of AI demands a comprehensive examination of its implications", $true, $false, $false, $false, $false, $true, 1, $false, "The diversity of cells is staggering, ranging from simple prokaryotes, like bacteria, to intricate eukaryotes, such as animal and plant cells", 2)

$d.Content.Find.Execute(" Our growing dependence on AI-powered systems in domains as diverse as healthcare, finance, and criminal justice mandates ethical scrutiny", $true, $false, $false, $false, $false, $true, 1, $false, " Each cell type is uniquely specialized, adapted to perform specific tasks essential for the survival of the organism", 2)

$d.Content.Find.Execute(" Are AI algorithms biased? Do they perpetuate existing prejudices? What are the consequences of AI-driven decisions gone awry? These are just a few of the ethical landmines we must navigate to ensure AI's ethical integrity", $true, $false, $false, $false, $false, $true, 1, $false, " Specialized cells, such as neurons, facilitate the rapid transmission of information throughout organisms, while muscle cells enable movement and contraction. This exquisite symphony of cells working in concert underscores the intricate complexity of life", 2)

$d.Content.Find.Execute("Furthermore, the rise of autonomous AI systems presents unprecedented challenges", $true, $false, $false, $false, $false, $true, 1, $false, "The study of cells has revolutionized our understanding of biology and medicine", 2)

$d.Content.Find.Execute(" As these machines become increasingly autonomous, the questions of accountability and liability become tangled", $true, $false, $false, $false, $false, $true, 1, $false, " The development of microscopes has allowed scientists to peer into the inner sanctums of cells, revealing the intricate structures and processes that govern life", 2)

$d.Content.Find.Execute(" Who bears responsibility when an AI system malfunctions or makes harmful decisions? Legal frameworks and ethical principles must evolve swiftly to address such conundrums", $true, $false, $false, $false, $false, $true, 1, $false, " This knowledge has led to breakthroughs in treating diseases, developing new drugs, and understanding the genetic basis of inheritance. The study of cells continues to unlock mysteries, pushing the boundaries of biological knowledge and offering hope for new treatments and therapies", 2)

# --- Summary paragraph ---
$d.Content.Find.Execute("AI's rapidly expanding role in our lives amplifies the need for ethical considerations", $true, $false, $false, $false, $false, $true, 1, $false, "The microscopic world of cells is a captivating realm of intricate structures and processes that hold the secrets of life", 2)

$d.Content.Find.Execute(" We must delve into the complexities of AI-driven decision-making, scrutinize potential biases, and contemplate the consequences of AI's ever-growing autonomy", $true, $false, $false, $false, $false, $true, 1, $false, " From the simplest prokaryotes to the complex eukaryotes, each cell is a finely tuned machine, performing specialized tasks essential for the survival of the organism", 2)

$d.Content.Find.Execute(" This ethical exploration encompasses concerns of privacy, transparency, accountability, and liability", $true, $false, $false, $false, $false, $true, 1, $false, " The study of cells has revolutionized biology and medicine, leading to groundbreaking discoveries that have improved our understanding of diseases, genetics, and treatments", 2)

$d.Content.Find.Execute(" Only by grappling with these challenges head-on can we harness AI's potential for progress while safeguarding our values and ensuring its ethical compass remains steadfast", $true, $false, $false, $false, $false, $true, 1, $false, " As we continue to explore the fascinating realm of cells, we unlock the mysteries of life and pave the way for new advancements in healthcare and biological knowledge", 2)

# --- Add a new empty paragraph at the very end of the document ---
$d.Paragraphs.Add()
